# Weekly price-sheet update:
# - A new week of data is inserted at row 67 (pushing existing rows down by one).
# - The oldest row (previously row 140) survives as the new last row (141).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 67 (shifts 67..140 down to 68..141).
$ws.Rows.Item(67).Insert()

# Populate the newly inserted row 67 with this week's record.
$ws.Range("A67").Value = 7
$ws.Range("B67").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C67").Value = "Ñuble"
$ws.Range("D67").Value = 45062
$ws.Range("E67").Value = 16
$ws.Range("F67").Value = 100112037
$ws.Range("G67").Value = "Cebollín"
$ws.Range("H67").Value = "Sin especificar"
$ws.Range("I67").Value = "Primera"
$ws.Range("J67").Value = 100
$ws.Range("K67").Value = 6000
$ws.Range("L67").Value = 6500
$ws.Range("M67").Value = 6250
$ws.Range("N67").Value = "$/paquete 36 unidades"
$ws.Range("O67").Value = "Provincia de Diguillín"
$ws.Range("P67").Value = 174
$ws.Range("Q67").Value = 36
$ws.Range("R67").Value = "Hortaliza"
